$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2020" column (N) -----------------------------------------
# N4: year label 2020, matching the style used by the other year-header cells
# (D4:K4, style index 12 in the original file).
$ws.Range("D4").Copy()
$ws.Range("N4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N4").Value = 2020

# N5: data value 534, matching the style used by the neighbouring value cells
# (L5:M5, style index 17 in the original file).
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N5").Value = 534

$excel.CutCopyMode = 0

# --- Update the view state --------------------------------------------------
# Scroll the window so column E is the left-most visible column, and select
# cell S10 (matches the sheetView/selection seen in the target workbook).
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("S10").Select()
